$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ENVELOPE_ASSEMBLIES")

# Update the basement (type_base, column I) assembly for STANDARD2..STANDARD6
# (rows 3-7) from FLOOR_AS4 to the new FLOOR_AS6, mirroring the existing
# floor (type_floor, column H) update.
$ws.Range("I3:I7").Value = "FLOOR_AS6"

# Reflect the new active selection left behind on this sheet.
$ws.Range("H17").Select()
